# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker table (rows 16-25, columns B:G) is re-sorted: originally the
# data was grouped by worker (all periods for Maria, then all periods for
# Martha); now it is grouped by period (both workers for period 1811, then
# both workers for 1812, etc). Column F (Valor Mora) keeps 26041 for period
# 1903 and 31249 for every other period.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$tipoDoc = "CC"
$ccMaria  = "33297383"
$nombreMaria  = "MARIA EUGENIA ROCHA PABUENA"
$ccMartha = "45528190"
$nombreMartha = "MARTHA LUCIA ROCHA PABUENA"
$salario = 781242

$periodos = @("1811", "1812", "1901", "1902", "1903")

$row = 16
foreach ($periodo in $periodos) {
    if ($periodo -eq "1903") { $valorMora = 26041 } else { $valorMora = 31249 }

    $ws.Cells.Item($row, 2).Value = $tipoDoc
    $ws.Cells.Item($row, 3).Value = $ccMaria
    $ws.Cells.Item($row, 4).Value = $nombreMaria
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $valorMora
    $ws.Cells.Item($row, 7).Value = $salario
    $row = $row + 1

    $ws.Cells.Item($row, 2).Value = $tipoDoc
    $ws.Cells.Item($row, 3).Value = $ccMartha
    $ws.Cells.Item($row, 4).Value = $nombreMartha
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $valorMora
    $ws.Cells.Item($row, 7).Value = $salario
    $row = $row + 1
}
